$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name and title text to reflect the new "through" date
$ws.Name = "Through 2021-11-22"

# Update the label for the November row (column A, row 12)
$ws.Range("A12").Value = "November (through 11-22)"

# Update October row (row 11) - only 2021 (column H) value changes
$ws.Range("H11").Value = 195

# Update November row (row 12) - all year columns change
$ws.Range("B12").Value = 22
$ws.Range("C12").Value = 53
$ws.Range("D12").Value = 88
$ws.Range("E12").Value = 43
$ws.Range("F12").Value = 36
$ws.Range("G12").Value = 152
$ws.Range("H12").Value = 151

# Update Total row (row 13) - all year columns change
$ws.Range("B13").Value = 280
$ws.Range("C13").Value = 539
$ws.Range("D13").Value = 798
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = 518
$ws.Range("G13").Value = 1209
$ws.Range("H13").Value = 1594
